$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.515.02'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '1.580.09'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = "'" + '208.23'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").Value = "'" + '22.32'
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = "'" + '0.249'
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("D12").Value = '1.809.50'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '1.575.67'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").Value = "'" + '3.83'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").Value = "'" + '0.525'
$ws.Range("E15").Value = '  -2.13%  '
$ws.Range("D16").Value = '27.553.21'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = "'" + '63.05'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = "'" + '214.48'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0691'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = "'" + '7.31'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").Value = "'" + '4.14'
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("D23").Value = "'" + '9.78'
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = "'" + '152.96'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = "'" + '6.93'
$ws.Range("E26").Value = '  +2.97%  '
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = "'" + '15.06'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("D32").Value = "'" + '3.22'
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("D33").Value = '1.371.70'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").Value = "'" + '1.55'
$ws.Range("E35").Value = '  +1.28%  '
$ws.Range("D36").Value = "'" + '0.971'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("E38").Value = '  +1.23%  '
$ws.Range("D39").Value = "'" + '0.533'
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("D40").Value = "'" + '0.826'
$ws.Range("E40").Value = '  +1.70%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  -0.45%  '
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("D44").Value = "'" + '64.44'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("E46").Value = '  +2.50%  '
$ws.Range("D47").Value = '1.720.37'
$ws.Range("D48").Value = "'" + '85.46'
$ws.Range("E48").Value = '  -2.77%  '
$ws.Range("D49").Value = '0.0₇0994'
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("E50").Value = '  -1.02%  '

Write-Host "Applied cryptos update"